# Update score lab 2 OOP
# Fills column D (lab 2 score) for the rows that previously had an empty
# inlineStr placeholder, using the values captured in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scores = @{
    3  = 600
    4  = 800
    8  = 785.71
    10 = 800
    11 = 800
    12 = 800
    13 = 800
    14 = 800
    15 = 300
    17 = 600
    18 = 800
    19 = 800
    20 = 600
    21 = 400
    22 = 500
    23 = 800
    25 = 600
    26 = 800
    29 = 800
    30 = 800
    31 = 800
    32 = 800
    34 = 800
    35 = 614.29
    37 = 800
    38 = 800
    39 = 800
    40 = 771.4299999999999
    41 = 800
    43 = 800
    44 = 771.4299999999999
    45 = 800
    46 = 800
    48 = 700
}

foreach ($row in $scores.Keys) {
    $ws.Range("D$row").Value = $scores[$row]
}
